$d = $word.ActiveDocument

$replacements = @(
    @("2024-08-19 Monday", "2024-08-20 Tuesday"),
    @("939÷6=156, 3", "232÷2=116, 0"),
    @("974÷6=162, 2", "435÷3=145, 0"),
    @("140÷8=17, 4", "920÷9=102, 2"),
    @("769÷5=153, 4", "206÷4=51, 2"),
    @("266÷6=44, 2", "987÷7=141, 0"),
    @("681÷6=113, 3", "149÷5=29, 4"),
    @("933÷9=103, 6", "637÷6=106, 1"),
    @("574÷6=95, 4", "125÷8=15, 5"),
    @("102÷3=34, 0", "857÷5=171, 2"),
    @("787÷6=131, 1", "461÷9=51, 2"),
    @("376÷4=94, 0", "665÷6=110, 5"),
    @("454÷4=113, 2", "251÷3=83, 2"),
    @("250÷6=41, 4", "436÷4=109, 0"),
    @("517÷2=258, 1", "232÷4=58, 0"),
    @("908÷4=227, 0", "398÷6=66, 2"),
    @("586÷2=293, 0", "323÷3=107, 2"),
    @("627÷2=313, 1", "887÷9=98, 5"),
    @("310÷4=77, 2", "655÷9=72, 7"),
    @("766÷2=383, 0", "122÷3=40, 2"),
    @("714÷9=79, 3", "956÷8=119, 4"),
    @("432÷6=72, 0", "367÷6=61, 1"),
    @("825÷5=165, 0", "456÷9=50, 6"),
    @("646÷2=323, 0", "380÷5=76, 0"),
    @("436÷3=145, 1", "963÷4=240, 3"),
    @("436÷5=87, 1", "671÷5=134, 1")
)

foreach ($pair in $replacements) {
    $find = $pair[0]
    $replace = $pair[1]
    $range = $d.Content
    $range.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
}
